$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 503
$ws.Range("I19").Value = 313.2857
$ws.Range("J19").Value = 558.3333
$ws.Range("K19").Value = 313.2857
$ws.Range("L19").Value = 558.3333
$ws.Range("M19").Value = -138.2857
$ws.Range("N19").Value = -908.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 63245.3
$ws.Range("I86").Value = 137778.44
$ws.Range("K86").Value = 137778.44
$ws.Range("M86").Value = -136655.44

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 63245.3
$ws.Range("I89").Value = 137778.44
$ws.Range("K89").Value = 688892.2
$ws.Range("M89").Value = -683276.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 4605.4585
$ws.Range("I111").Value = 2884.2144
$ws.Range("J111").Value = 7015.2
$ws.Range("K111").Value = 8652.643199999999
$ws.Range("L111").Value = 21045.6
$ws.Range("M111").Value = -5585.643199999999
$ws.Range("N111").Value = -27179.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 8930279
$ws.Range("I137").Value = 15153389
$ws.Range("J137").Value = 1469.0435
$ws.Range("K137").Value = 45460167
$ws.Range("L137").Value = 4407.1305
$ws.Range("M137").Value = -45457617
$ws.Range("N137").Value = -9507.130499999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2374.7874
$ws.Range("J138").Value = 2533.3333
$ws.Range("L138").Value = 7599.999899999999
$ws.Range("N138").Value = -17879.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4913993.5
$ws.Range("I32").Value = 6774.6895
$ws.Range("K32").Value = 6774.6895
$ws.Range("M32").Value = -6487.6895

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4116.2666
$ws.Range("I45").Value = 2811.8333
$ws.Range("J45").Value = 4985.8887
$ws.Range("K45").Value = 2811.8333
$ws.Range("L45").Value = 4985.8887
$ws.Range("M45").Value = -2434.8333
$ws.Range("N45").Value = -5739.8887

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 30613146
$ws.Range("I74").Value = 37500856
$ws.Range("J74").Value = 1108
$ws.Range("K74").Value = 37500856
$ws.Range("L74").Value = 1108
$ws.Range("M74").Value = -37499982
$ws.Range("N74").Value = -2856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 30613146
$ws.Range("I77").Value = 37500856
$ws.Range("J77").Value = 1108
$ws.Range("K77").Value = 187504280
$ws.Range("L77").Value = 5540
$ws.Range("M77").Value = -187499912
$ws.Range("N77").Value = -14276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2904.1538
$ws.Range("I122").Value = 1330.4445
$ws.Range("K122").Value = 3991.3335
$ws.Range("M122").Value = -1541.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2362275.2
$ws.Range("I132").Value = 1853.12
$ws.Range("J132").Value = 4469795
$ws.Range("K132").Value = 5559.36
$ws.Range("L132").Value = 13409385
$ws.Range("M132").Value = -3029.36
$ws.Range("N132").Value = -13414445

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1586.6666
$ws.Range("I86").Value = 1406
$ws.Range("K86").Value = 1406
$ws.Range("M86").Value = -283

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1586.6666
$ws.Range("I89").Value = 1406
$ws.Range("K89").Value = 7030
$ws.Range("M89").Value = -1414

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2799.1
$ws.Range("I107").Value = 2788.875
$ws.Range("J107").Value = 2840
$ws.Range("K107").Value = 2788.875
$ws.Range("L107").Value = 2840
$ws.Range("M107").Value = -868.875
$ws.Range("N107").Value = -6680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 71431944
$ws.Range("I122").Value = 125001400
$ws.Range("J122").Value = 6002.5
$ws.Range("K122").Value = 375004200
$ws.Range("L122").Value = 18007.5
$ws.Range("M122").Value = -375001750
$ws.Range("N122").Value = -22907.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2581.625
$ws.Range("I132").Value = 1800
$ws.Range("K132").Value = 5400
$ws.Range("M132").Value = -2870

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 455414.62
$ws.Range("I68").Value = 612.04877
$ws.Range("J68").Value = 771463.9
$ws.Range("K68").Value = 1836.14631
$ws.Range("L68").Value = 2314391.7
$ws.Range("M68").Value = -1025.14631
$ws.Range("N68").Value = -2316013.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 455414.62
$ws.Range("I71").Value = 612.04877
$ws.Range("J71").Value = 771463.9
$ws.Range("K71").Value = 5508.43893
$ws.Range("L71").Value = 6943175.100000001
$ws.Range("M71").Value = -1452.43893
$ws.Range("N71").Value = -6951287.100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 26871368
$ws.Range("I107").Value = 33335112
$ws.Range("J107").Value = 2632328.5
$ws.Range("K107").Value = 100005336
$ws.Range("L107").Value = 7896985.5
$ws.Range("M107").Value = -100003416
$ws.Range("N107").Value = -7900825.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 471.92856
$ws.Range("I113").Value = 437.64285
$ws.Range("K113").Value = 1312.92855
$ws.Range("M113").Value = 857.0714499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 769.0741
$ws.Range("I107").Value = 398.75
$ws.Range("J107").Value = 3731.6667
$ws.Range("K107").Value = 398.75
$ws.Range("L107").Value = 3731.6667
$ws.Range("M107").Value = 1521.25
$ws.Range("N107").Value = -7571.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 716301.0600000001
$ws.Range("I122").Value = 1112289.6
$ws.Range("J122").Value = 3521.6
$ws.Range("K122").Value = 3336868.8
$ws.Range("L122").Value = 10564.8
$ws.Range("M122").Value = -3334418.8
$ws.Range("N122").Value = -15464.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 45000
$ws.Range("J138").Value = 45000
$ws.Range("L138").Value = 45000
$ws.Range("N138").Value = -55280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 50502230
$ws.Range("I40").Value = 63126176
$ws.Range("J40").Value = 6451.25
$ws.Range("K40").Value = 63126176
$ws.Range("L40").Value = 6451.25
$ws.Range("M40").Value = -63126040
$ws.Range("N40").Value = -6723.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 36250
$ws.Range("J122").Value = 17500
$ws.Range("L122").Value = 52500
$ws.Range("N122").Value = -57400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 16500
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 16500
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 16500
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -17084

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9012.666999999999
$ws.Range("I122").Value = 11071.857
$ws.Range("K122").Value = 33215.571
$ws.Range("M122").Value = -30765.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 48999.8
$ws.Range("J123").Value = 48999.8
$ws.Range("L123").Value = 48999.8
$ws.Range("N123").Value = -58799.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6711.769
$ws.Range("I126").Value = 7182.4346
$ws.Range("J126").Value = 3103.3333
$ws.Range("K126").Value = 21547.3038
$ws.Range("L126").Value = 9309.999899999999
$ws.Range("M126").Value = -19077.3038
$ws.Range("N126").Value = -14249.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8073940
$ws.Range("I136").Value = 9625958
$ws.Range("J136").Value = 3445
$ws.Range("K136").Value = 28877874
$ws.Range("L136").Value = 10335
$ws.Range("M136").Value = -28875324
$ws.Range("N136").Value = -15435
